$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- sheet1 ("data"): the request row now only keeps the manager name + full request text ---
$ws1.Range("B2").ClearContents()
$ws1.Range("C2").ClearContents()
$ws1.Range("E2").ClearContents()
$ws1.Range("F2").ClearContents()

$appText = @'
\nИГО:\nЗаявка на доставку \n1. Дата отгрузки 28.04.2023 \n2. Марка ЦЕМ I 42.5н Беларусь   \n3. Количество тонн: 100 \n4. От ООО Спарта \n5. Завод: Сзтк \n6. Покупатель ООО ""ТД"Цемент \n7. Грузополучатель: ООО "ТД"Цемент  \n8. Голицыно\n+7 910 404-06-14\nРБУ\nМожайское ш., 81\n
'@
$ws1.Range("A2").Value = $appText

# --- add the new "errors" sheet right after "data" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "ошибки"

# headers reuse the same bold / bordered style as sheet1's header row
$ws1.Range("A1").Copy()
$ws2.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("A1").Value = "Ошибка"
$ws2.Range("B1").Value = "Заявка"

$traceback = @'
Traceback (most recent call last):
  File "C:\Users\79852\Desktop\gh\applications\applic.py", line 64, in find_unit_note
    unit_found = self.units_df.loc[self.units_df[self.units_df['data'] == unit].index[0], 'unit']
  File "C:\Users\79852\anaconda3\lib\site-packages\pandas\core\indexes\base.py", line 5039, in __getitem__
    return getitem(key)
IndexError: index 0 is out of bounds for axis 0 with size 0
During handling of the above exception, another exception occurred:
Traceback (most recent call last):
  File "C:\Users\79852\Desktop\gh\applications\applic.py", line 91, in process_application
    unit_found = self.find_unit_note(application)
  File "C:\Users\79852\Desktop\gh\applications\applic.py", line 67, in find_unit_note
    self.error_log.append(f"Ошибка в методе find_unit_note: {e}")
  File "C:\Users\79852\anaconda3\lib\site-packages\pandas\core\frame.py", line 9039, in append
    return self._append(other, ignore_index, verify_integrity, sort)
  File "C:\Users\79852\anaconda3\lib\site-packages\pandas\core\frame.py", line 9082, in _append
    result = concat(
  File "C:\Users\79852\anaconda3\lib\site-packages\pandas\util\_decorators.py", line 311, in wrapper
    return func(*args, **kwargs)
  File "C:\Users\79852\anaconda3\lib\site-packages\pandas\core\reshape\concat.py", line 347, in concat
    op = _Concatenator(
  File "C:\Users\79852\anaconda3\lib\site-packages\pandas\core\reshape\concat.py", line 437, in __init__
    raise TypeError(msg)
TypeError: cannot concatenate object of type '<class 'str'>'; only Series and DataFrame objs are valid

'@
$ws2.Range("A2").Value = $traceback
$ws2.Range("B2").Value = $appText

$ws2.Range("A2").WrapText = $true
$ws2.Rows.Item(2).RowHeight = 96.5
$ws2.Columns.Item(1).ColumnWidth = 46.333333333333336

$ws2.Range("C2").Select()
$ws1.Activate()
